# Students_Details.xlsx: add Age / Salary / Department columns (M:O)
# to the student table on Sheet1, driven from a headers-based POJO
# (web-table "Add new record" refactor).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New header row (row 1) --------------------------------------------
$ws.Range("M1").Value = "Age"
$ws.Range("N1").Value = "Salary"
$ws.Range("O1").Value = "Department"

# --- New per-record data (rows 2-5), matching existing row order -------
$newData = @(
    @{ Row = 2; Age = 34; Salary = 28000; Department = "Information Technology" },
    @{ Row = 3; Age = 30; Salary = 18000; Department = "Finance" },
    @{ Row = 4; Age = 40; Salary = 15000; Department = "Medicine" },
    @{ Row = 5; Age = 42; Salary = 19000; Department = "Recreation" }
)

foreach ($rec in $newData) {
    $r = $rec.Row
    $ws.Cells.Item($r, 13).Value = $rec.Age          # M
    $ws.Cells.Item($r, 14).Value = $rec.Salary       # N
    $ws.Cells.Item($r, 15).Value = $rec.Department   # O
}

# --- Restore a sensible selection / scroll position over the new cols --
$ws.Range("M6").Select()
$excel.ActiveWindow.ScrollColumn = 3
